$d = $word.ActiveDocument

# The paragraph that introduces the list of orphan tags.
$introPara = $d.Paragraphs.Item(2)
$introRng = $introPara.Range

# Insert a fresh (initially empty) paragraph right after it; this will
# serve as the running insertion point for the new list entries.
[void]$introRng.InsertParagraphAfter()

# One line per orphan tag found in the documents, in the same order as
# the original report.
$orphanTags = @('PUMP:RISK:10 ', 'PUMP:RISK:20 ', 'PUMP:RISK:30 ', 'PUMP:RISK:40 ', 'PUMP:RISK:50 ', 'PUMP:URS:1 ', 'PUMP:URS:3 ', 'PUMP:URS:8 ', 'PUMP:URS:10 ', 'PUMP:URS:100 ', 'PUMP:URS:103 ', 'PUMP:URS:1000 ', 'PUMP:URS:3330 ', 'PUMP:URS:3350 ', 'PUMP:URS:4000 ', 'PUMP:HRS:103', 'PUMP:TBV:1111', 'PUMP:PRS:103', 'ACE:SRS:110', 'ACE:SRS:120', 'PUMP:TBV:1', 'PUMP:PRS:6', 'PUMP:TBD:1', 'PUMP:DER:2', 'ACE:SRS:1000', 'PUMP:UNIT:100', 'PUMP:UNIT:110', 'PUMP:UNIT:120', 'PUMP:UNIT:130', 'PUMP:UNIT:140', 'PUMP:UNIT:150', 'PUMP:UNIT:160', 'PUMP:UNIT:170', 'PUMP:UNIT:180', 'PUMP:UNIT:190', 'PUMP:UNIT:200', 'PUMP:UNIT:210', 'PUMP:UNIT:220')

foreach ($tag in $orphanTags) {
    # Always (re)target the current last paragraph of the document: each
    # InsertXML call below adds the new paragraph before this carrier
    # paragraph mark, so re-fetching it keeps the insertion point moving
    # forward one entry at a time.
    $carrierPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $carrierRng = $carrierPara.Range

    $escaped = $tag -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    # Only mark the text run as whitespace-preserving when the tag
    # actually has leading/trailing whitespace that needs protecting,
    # matching how Word itself emits xml:space="preserve".
    if ($tag -ne $tag.Trim()) {
        $tTag = '<w:t xml:space="preserve">' + $escaped + '</w:t>'
    } else {
        $tTag = '<w:t>' + $escaped + '</w:t>'
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $tTag + '</w:r></w:p></w:body>' +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$carrierRng.InsertXML($xml)
}

# The loop above leaves one stray empty paragraph at the very end of the
# document (the original carrier paragraph mark, now empty). Merge it
# into the preceding paragraph by deleting the span from the end of the
# previous paragraph's content through the end of this empty one.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
if (($lastPara.Range.End - $lastPara.Range.Start) -le 1) {
    $mergeRng = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
    [void]$mergeRng.Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
